# Fixed deployment dates based on cruise reports and WHOI documentation
# Global Irminger - GI05MOAS-GL485 deployment info

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Moorings")

# Recover Date (G2) - set to 2015-11-20 (serial 42328)
$ws.Range("G2").Value = Get-Date -Year 2015 -Month 11 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Notes (L2) - glider lost
$ws.Range("L2").Value = "glider lost"

# Update the active selection on the sheet to E11 (as last edited cell)
$ws.Activate()
$ws.Range("E11").Select()
